$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 36; this shifts existing rows 36-39 down to 37-40
$ws.Rows.Item(36).Insert()

# Populate the new row 36 with the new data entry
$ws.Cells.Item(36, 1).Value = 11
$ws.Cells.Item(36, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(36, 3).Value = "Bíobío"
$ws.Cells.Item(36, 4).Value = 44615
$ws.Cells.Item(36, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(36, 5).Value = 8
$ws.Cells.Item(36, 6).Value = 100112043
$ws.Cells.Item(36, 7).Value = "Pepino dulce"
$ws.Cells.Item(36, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 100
$ws.Cells.Item(36, 11).Value = 15000
$ws.Cells.Item(36, 12).Value = 16000
$ws.Cells.Item(36, 13).Value = 15500
$ws.Cells.Item(36, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(36, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(36, 16).Value = 861
$ws.Cells.Item(36, 17).Value = 18
$ws.Cells.Item(36, 18).Value = "Hortaliza"

$ws.Range("A1").Select()
